# Generate Report for Handoff
# Replaces the pair of e2e test files tracked in the localization status
# report with a new pair, and updates their handoff/handback status.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "53f9b59a-0d00-4773-9a77-bffcd1475d57"
$newGuid1 = "c4ebf2e4-7658-4fa0-8ad3-b2819fe58533"
$oldGuid2 = "822d4394-994b-4792-8283-39888d3d852c"
$newGuid2 = "ffff60a0a25e-bc68-4bc8-b309-195936628f89"

$newFile1 = "$newGuid1.md"
$newFile2 = "$newGuid2.md"
$newPath1 = "e2e\$newFile1"
$newPath2 = "e2e\$newFile2"

$newStatus = "Ready for handoff"
$newOverviewDate = "2016-08-12 13:15:11"

$newZhXlf = "$newGuid1.480074fc3d12cdf0b233289320956d97e9880ccf.zh-cn.xlf"
$newDeXlf = "$newGuid1.480074fc3d12cdf0b233289320956d97e9880ccf.de-de.xlf"
$newZhHandoffDate = "2016-08-12 13:14:57"
$newDeHandoffDate = "2016-08-12 13:15:11"
$newHandbackDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $newPath1
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newOverviewDate

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $newPath2
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $newOverviewDate

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = $newPath1
    } elseif ($addr -eq '$B$3') {
        $h.TextToDisplay = $newPath2
    }
}

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("K2").Value = $newHandbackDate

$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("H3").Value = $newZhHandoffDate
$wsZh.Range("K3").Value = $newHandbackDate

$zhToDelete = @()
foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newFile1
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $newFile2
    } elseif ($addr -eq '$I$2' -or $addr -eq '$I$3') {
        $zhToDelete += $h
    }
}
foreach ($h in $zhToDelete) {
    $h.Delete()
}
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("I3").Value = ""
$wsZh.Range("I3").Style = "Normal"
$wsZh.Range("J3").Value = ""

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newDeHandoffDate
$wsDe.Range("K2").Value = $newHandbackDate

$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("H3").Value = $newDeHandoffDate
$wsDe.Range("K3").Value = $newHandbackDate

$deToDelete = @()
foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newFile1
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $newFile2
    } elseif ($addr -eq '$I$2' -or $addr -eq '$I$3') {
        $deToDelete += $h
    }
}
foreach ($h in $deToDelete) {
    $h.Delete()
}
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("I3").Value = ""
$wsDe.Range("I3").Style = "Normal"
$wsDe.Range("J3").Value = ""
